$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain text, matching the
# scraped data format (values like "1.00" or "6.96" must not be coerced
# into numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '60.875.70'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '3.413.23'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '572.16'
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").Value = '141.87'
$ws.Range("E6").Value = '  -4.15%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.409.99'
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '0.478'
$ws.Range("E9").Value = '  +0.77%  '
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '3.998.69'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = '3.418.49'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '61.020.99'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("D19").Value = '6.33'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").Value = '14.39'
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").Value = '9.34'
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("D22").Value = '392.69'
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '73.02'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").Value = '0.995'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").Value = '3.568.44'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").Value = '  -3.85%  '
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").Value = '8.15'
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("D33").Value = '1.45'
$ws.Range("E33").Value = '  -6.25%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '23.85'
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("D36").Value = '6.99'
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D37").Value = '3.441.79'
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("D38").Value = '5.13'
$ws.Range("E38").Value = '  -2.70%  '
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("D40").Value = '167.13'
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = '27.07'
$ws.Range("E42").Value = '  +4.11%  '
$ws.Range("D43").Value = '0.796'
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  +1.03%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.71'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = '41.81'
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").Value = '2.606.55'
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("E49").Value = '  -4.72%  '
$ws.Range("D50").Value = '6.95'
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("D51").Value = '22.95'
$ws.Range("E51").Value = '  -4.03%  '
